$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 7 - Restore Running (Maintenance Tasks Running)
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Maintenance Tasks Running"
$ws.Range("D7").Value = "Restore Running"
$ws.Hyperlinks.Add($ws.Range("E7"), "http://BrentOzar.com/go/backups")

# Row 8 - Data File Growing (SQL Server Internal Maintenance)
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "SQL Server Internal Maintenance"
$ws.Range("D8").Value = "Data File Growing"
$ws.Hyperlinks.Add($ws.Range("E8"), "http://BrentOzar.com/go/ifi")

# Row 9 - Long-Running Query Blocking Others (Query Problems)
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Query Problems"
$ws.Range("D9").Value = "Long-Running Query Blocking Others"
$ws.Hyperlinks.Add($ws.Range("E9"), "http://BrentOzar.com/go/blocking")

$ws.Range("A10").Select()
